# Updates cryptos list values per the latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.716.82"
$ws.Range("E2").Value = "  -1.47%  "
$ws.Range("D3").Value = "'3.386.00"
$ws.Range("E3").Value = "  -1.89%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'567.66"
$ws.Range("E5").Value = "  -2.32%  "
$ws.Range("D6").Value = "'141.24"
$ws.Range("E6").Value = "  -3.23%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'3.386.34"
$ws.Range("E8").Value = "  -1.91%  "
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("D11").Value = "'0.124"
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("D12").Value = "'0.397"
$ws.Range("E12").Value = "  +1.68%  "
$ws.Range("D13").Value = "'3.962.30"
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("D14").Value = "'28.31"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("E15").Value = "  +1.85%  "
$ws.Range("D16").Value = "'0.0000171"
$ws.Range("E16").Value = "  -2.10%  "
$ws.Range("D17").Value = "'3.379.70"
$ws.Range("E17").Value = "  -2.16%  "
$ws.Range("D18").Value = "'60.830.75"
$ws.Range("E18").Value = "  -1.48%  "
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Value = "'14.03"
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("D21").Value = "'9.02"
$ws.Range("E21").Value = "  -5.48%  "
$ws.Range("D22").Value = "'383.50"
$ws.Range("E22").Value = "  -1.71%  "
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("D24").Value = "'73.68"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("E26").Value = "  -5.40%  "
$ws.Range("D27").Value = "'3.522.22"
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("E28").Value = "  -1.40%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").Value = "'7.41"
$ws.Range("E30").Value = "  -2.98%  "
$ws.Range("E31").Value = "  -2.42%  "
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("E33").Value = "  -3.24%  "
$ws.Range("D35").Value = "'23.72"
$ws.Range("E35").Value = "  -1.56%  "
$ws.Range("E36").Value = "  -0.48%  "
$ws.Range("D37").Value = "'166.44"
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("E38").Value = "  -2.36%  "
$ws.Range("D39").Value = "'3.415.44"
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("E40").Value = "  -4.78%  "
$ws.Range("D41").Value = "'27.91"
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("D42").Value = "'0.0776"
$ws.Range("E42").Value = "  -1.33%  "
$ws.Range("E43").Value = "  -3.06%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "'41.83"
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("E46").Value = "  -2.37%  "
$ws.Range("E47").Value = "  -2.99%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'1.12"
$ws.Range("E48").Value = "  -3.34%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'2.518.25"
$ws.Range("E49").Value = "  -2.01%  "
$ws.Range("D50").Value = "'23.68"
$ws.Range("E50").Value = "  +2.68%  "
$ws.Range("D51").Value = "'6.83"
$ws.Range("E51").Value = "  -2.02%  "
